$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns being updated so that numeric-looking
# strings (e.g. "300.93", "-0.92%") are preserved as text, matching the source data.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "E20", "D21", "E21", "D22", "E22", "D24", "E24", "D25", "E25", "E26", "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "E43", "D44", "E44", "D45", "E45", "E46", "E47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "300.93"
$ws.Range("E2").Value = "-0.92%"
$ws.Range("D3").Value = "31.36"
$ws.Range("E3").Value = "-2.64%"
$ws.Range("D4").Value = "5.150"
$ws.Range("E4").Value = "-2.35%"
$ws.Range("D5").Value = "0.07371"
$ws.Range("E5").Value = "-1.52%"
$ws.Range("D6").Value = "2.432"
$ws.Range("E6").Value = "59.93%"
$ws.Range("D7").Value = "7.952"
$ws.Range("E7").Value = "1.23%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "3.789"
$ws.Range("E8").Value = "-0.52%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9206"
$ws.Range("E9").Value = "0.03%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "0.1730"
$ws.Range("E10").Value = "2.67%"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "0.07646"
$ws.Range("E11").Value = "-4.63%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08097"
$ws.Range("E12").Value = "0.53%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03027"
$ws.Range("E13").Value = "0.65%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09929"
$ws.Range("E14").Value = "0.24%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001501"
$ws.Range("E15").Value = "0.22%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006134"
$ws.Range("E16").Value = "-4.92%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.465"
$ws.Range("E17").Value = "-0.05%"
$ws.Range("E18").Value = "-0.22%"
$ws.Range("E20").Value = "-0.45%"
$ws.Range("D21").Value = "4.652"
$ws.Range("E21").Value = "3.73%"
$ws.Range("D22").Value = "0.04653"
$ws.Range("E22").Value = "0.92%"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").Value = "0.60%"
$ws.Range("D25").Value = "0.004487"
$ws.Range("E25").Value = "0.82%"
$ws.Range("E26").Value = "-7.13%"
$ws.Range("E27").Value = "5.43%"
$ws.Range("D39").Value = "0.01729"
$ws.Range("E39").Value = "0.41%"
$ws.Range("D40").Value = "0.04524"
$ws.Range("E40").Value = "0.62%"
$ws.Range("D41").Value = "0.007161"
$ws.Range("E41").Value = "0.08%"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").Value = "-0.09%"
$ws.Range("E43").Value = "-0.43%"
$ws.Range("D44").Value = "0.01072"
$ws.Range("E44").Value = "-16.24%"
$ws.Range("D45").Value = "0.00006273"
$ws.Range("E45").Value = "1.73%"
$ws.Range("E46").Value = "-22.98%"
$ws.Range("E47").Value = "171.71%"
